$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row => D value, E value ($null means "no change for that cell")
$updates = @{
    2  = @("28.250.82", "  +4.09%  ")
    3  = @("1.785.87", "  +0.17%  ")
    4  = @("0.9986", "  -0.53%  ")
    5  = @("336.93", "  +0.12%  ")
    6  = @("0.9961", $null)
    7  = @("0.3834", "  +0.31%  ")
    8  = @("0.3445", "  +0.96%  ")
    9  = @("47.62", "  -0.70%  ")
    10 = @("1.159", "  -2.27%  ")
    11 = @("0.07442", "  -0.01%  ")
    12 = @("23.20", "  +7.18%  ")
    13 = @("0.9951", "  -0.62%  ")
    14 = @("6.431", "  +0.03%  ")
    15 = @("1.784.03", "  +0.11%  ")
    16 = @("7.137", $null)
    17 = @("0.00001085", "  -0.57%  ")
    18 = @("0.06666", "  +0.40%  ")
    19 = @("82.89", "  -0.61%  ")
    20 = @("0.9960", "  -0.45%  ")
    21 = @("17.55", "  +1.08%  ")
    22 = @("6.446", "  -1.04%  ")
    23 = @("28.240.93", "  +4.04%  ")
    24 = @("12.14", "  -0.89%  ")
    25 = @("2.378", "  -0.10%  ")
    26 = @("20.98", "  -0.58%  ")
    27 = @("1.438", "  -0.23%  ")
    28 = @("2.423", "  -2.84%  ")
    29 = @("154.46", "  -0.21%  ")
    30 = @("1.985.72", "  +0.16%  ")
    31 = @("135.24", "  +0.95%  ")
    32 = @("6.188", "  +2.64%  ")
    33 = @($null, "  -0.76%  ")
    34 = @("0.08822", "  +1.99%  ")
    35 = @("12.82", "  -1.03%  ")
    36 = @("0.02436", "  +4.78%  ")
    37 = @("0.6897", "  +1.10%  ")
    38 = @("5.351", "  -0.58%  ")
    39 = @("0.06364", "  +1.11%  ")
    40 = @("0.2188", "  +0.61%  ")
    41 = @("1.246", "  +0.14%  ")
    42 = @("1.505", "  -7.42%  ")
    43 = @("8.356", "  -0.02%  ")
    44 = @("14.26", "  +0.59%  ")
    45 = @("0.9960", "  -0.36%  ")
    46 = @("0.6344", "  -0.95%  ")
    47 = @("3.854", "  +0.08%  ")
    48 = @("132.38", "  +0.82%  ")
    49 = @("2.105", "  -1.21%  ")
    50 = @("0.07478", "  +5.35%  ")
    51 = @("1.298", "  +10.44%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    if ($null -ne $dVal) {
        # Force the cell to remain text (these "Price" values are plain
        # strings in the source data, e.g. "28.250.82" / "0.9986"), then
        # restore the default "Normal" style so no formatting is left behind.
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.Style = "Normal"
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

# Row 51: coin renamed from EOS to Stacks, with a new coinranking.com link
$ws.Cells.Item(51, 2).Value = "Stacks"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
